$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency list on the worksheet (hourly price/volume refresh, GitHub Actions).
# Cells whose new value is a plain parseable number (e.g. "1.00", "0.0939") are given an
# explicit text NumberFormat first, so Excel keeps the exact digits/trailing zeros instead of
# silently converting the cell to a numeric value.

# Row 2
$ws.Cells.Item(2, 4).Value = '54.237.80'
$ws.Cells.Item(2, 5).Value = '  -2.75%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.284.85'
$ws.Cells.Item(3, 5).Value = '  -2.50%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '493.59'
$ws.Cells.Item(5, 5).Value = '  -2.08%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '127.00'
$ws.Cells.Item(6, 5).Value = '  -1.56%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.22%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -1.94%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.284.51'
$ws.Cells.Item(9, 5).Value = '  -2.85%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0939'
$ws.Cells.Item(10, 5).Value = '  -3.34%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.31%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.320'
$ws.Cells.Item(12, 5).Value = '  +0.14%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.62'
$ws.Cells.Item(13, 5).Value = '  -3.12%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '2.681.93'
$ws.Cells.Item(14, 5).Value = '  -2.83%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '21.51'
$ws.Cells.Item(15, 5).Value = '  -0.65%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '54.109.84'
$ws.Cells.Item(16, 5).Value = '  -2.83%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.271.29'
$ws.Cells.Item(18, 5).Value = '  -4.00%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.36%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +1.18%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '301.54'
$ws.Cells.Item(21, 5).Value = '  -2.90%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.40'
$ws.Cells.Item(22, 5).Value = '  +2.61%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.999'
$ws.Cells.Item(23, 5).Value = '  +0.11%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -2.55%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '63.69'
$ws.Cells.Item(25, 5).Value = '  -2.41%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.29%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.41%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '2.368.54'
$ws.Cells.Item(28, 5).Value = '  -3.34%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +1.55%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.43%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '165.22'
$ws.Cells.Item(31, 5).Value = '  -3.44%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.60'
$ws.Cells.Item(32, 5).Value = '  -2.52%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -3.05%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +1.79%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +0.01%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.999'
$ws.Cells.Item(36, 5).Value = '  +0.21%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +0.47%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '17.56'
$ws.Cells.Item(38, 5).Value = '  -0.63%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +1.23%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +5.41%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.38%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '35.39'
$ws.Cells.Item(42, 5).Value = '  -1.80%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +1.17%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +1.11%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.04%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'RenderToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '4.79'
$ws.Cells.Item(46, 5).Value = '  -0.48%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '125.61'
$ws.Cells.Item(47, 5).Value = '  -0.78%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0886'
$ws.Cells.Item(48, 5).Value = '  -0.45%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -1.94%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '236.69'
$ws.Cells.Item(50, 5).Value = '  -0.50%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0479'
$ws.Cells.Item(51, 5).Value = '  +1.00%  '
